$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 37740
$ws.Range("D2").Value = 54581015
$ws.Range("C3").Value = 90984
$ws.Range("D3").Value = 133372563
$ws.Range("C4").Value = 31180
$ws.Range("D4").Value = 46176349
$ws.Range("C5").Value = 8694
$ws.Range("D5").Value = 12921563
$ws.Range("C6").Value = 1995
$ws.Range("D6").Value = 2965006
$ws.Range("C12").Value = 41332
$ws.Range("D12").Value = 56077572
$ws.Range("C13").Value = 9651
$ws.Range("D13").Value = 13958458
$ws.Range("C14").Value = 25947
$ws.Range("D14").Value = 38053616
$ws.Range("C16").Value = 2152
$ws.Range("D16").Value = 3200165
$ws.Range("C20").Value = 10224
$ws.Range("D20").Value = 13537263
$ws.Range("C21").Value = 13379
$ws.Range("D21").Value = 19318492
$ws.Range("C22").Value = 31656
$ws.Range("D22").Value = 46454820
$ws.Range("C23").Value = 10219
$ws.Range("D23").Value = 15191178
$ws.Range("C24").Value = 2639
$ws.Range("D24").Value = 3923682
$ws.Range("C27").Value = 11686
$ws.Range("D27").Value = 15609553
$ws.Range("C28").Value = 7642
$ws.Range("D28").Value = 11069617
$ws.Range("C29").Value = 22480
$ws.Range("D29").Value = 32998049
$ws.Range("C30").Value = 7815
$ws.Range("D30").Value = 11630133
$ws.Range("C31").Value = 1959
$ws.Range("D31").Value = 2922999
$ws.Range("C34").Value = 8310
$ws.Range("D34").Value = 10976394
$ws.Range("C36").Value = 7828
$ws.Range("D36").Value = 11431934
$ws.Range("C37").Value = 3178
$ws.Range("D37").Value = 4709961
$ws.Range("C41").Value = 2472
$ws.Range("D41").Value = 3341353
$ws.Range("C42").Value = 17234
$ws.Range("D42").Value = 24918378
$ws.Range("C43").Value = 51104
$ws.Range("D43").Value = 74919486
$ws.Range("C44").Value = 19014
$ws.Range("D44").Value = 28243443
$ws.Range("C45").Value = 5605
$ws.Range("D45").Value = 8346677
$ws.Range("C50").Value = 16700
$ws.Range("D50").Value = 22233065
$ws.Range("C51").Value = 2022
$ws.Range("D51").Value = 2932471
$ws.Range("C52").Value = 6901
$ws.Range("D52").Value = 10144579
$ws.Range("C57").Value = 6985
$ws.Range("D57").Value = 9601687
$ws.Range("C58").Value = 946
$ws.Range("D58").Value = 1388579
$ws.Range("C59").Value = 2381
$ws.Range("D59").Value = 3530337
$ws.Range("C60").Value = 945
$ws.Range("D60").Value = 1407001
$ws.Range("C64").Value = 1393
$ws.Range("D64").Value = 1960706
$ws.Range("C65").Value = 15357
$ws.Range("D65").Value = 22181806
$ws.Range("C66").Value = 44673
$ws.Range("D66").Value = 65373153
$ws.Range("C67").Value = 15699
$ws.Range("D67").Value = 23330686
$ws.Range("C68").Value = 4568
$ws.Range("D68").Value = 6804292
$ws.Range("C69").Value = 924
$ws.Range("D69").Value = 1374168
$ws.Range("C73").Value = 15082
$ws.Range("D73").Value = 19884184
$ws.Range("C74").Value = 51382
$ws.Range("D74").Value = 74770492
$ws.Range("C75").Value = 146036
$ws.Range("D75").Value = 215143083
$ws.Range("C76").Value = 63622
$ws.Range("D76").Value = 94805657
$ws.Range("C77").Value = 20337
$ws.Range("D77").Value = 30385831
$ws.Range("C78").Value = 4816
$ws.Range("D78").Value = 7193043
$ws.Range("C85").Value = 50808
$ws.Range("D85").Value = 69112413
$ws.Range("C86").Value = 4602
$ws.Range("D86").Value = 6667513
$ws.Range("C87").Value = 11563
$ws.Range("D87").Value = 16987319
$ws.Range("C88").Value = 3883
$ws.Range("D88").Value = 5787083
$ws.Range("C89").Value = 1344
$ws.Range("D89").Value = 2008489
$ws.Range("C93").Value = 5411
$ws.Range("D93").Value = 7274275
$ws.Range("C94").Value = 1596
$ws.Range("D94").Value = 2298932
$ws.Range("C95").Value = 5164
$ws.Range("D95").Value = 7604743
$ws.Range("C101").Value = 3562
$ws.Range("D101").Value = 4715264
$ws.Range("C102").Value = 602
$ws.Range("D102").Value = 896664
$ws.Range("C107").Value = 10750
$ws.Range("D107").Value = 15594962
$ws.Range("C108").Value = 29200
$ws.Range("D108").Value = 42901815
$ws.Range("C109").Value = 9774
$ws.Range("D109").Value = 14534650
$ws.Range("C110").Value = 2685
$ws.Range("D110").Value = 4003707
$ws.Range("C114").Value = 9792
$ws.Range("D114").Value = 12934825
$ws.Range("C115").Value = 30447
$ws.Range("D115").Value = 43903358
$ws.Range("C116").Value = 66127
$ws.Range("D116").Value = 96774177
$ws.Range("C117").Value = 21363
$ws.Range("D117").Value = 31748940
$ws.Range("C124").Value = 25841
$ws.Range("D124").Value = 34513855
$ws.Range("C125").Value = 35971
$ws.Range("D125").Value = 51914426
$ws.Range("C126").Value = 76763
$ws.Range("D126").Value = 112249269
$ws.Range("C127").Value = 23839
$ws.Range("D127").Value = 35379909
$ws.Range("C129").Value = 1236
$ws.Range("D129").Value = 1838411
$ws.Range("C133").Value = 31808
$ws.Range("D133").Value = 42236713
$ws.Range("C134").Value = 13213
$ws.Range("D134").Value = 19125243
$ws.Range("C135").Value = 32315
$ws.Range("D135").Value = 47463152
$ws.Range("C137").Value = 2957
$ws.Range("D137").Value = 4408214
$ws.Range("C138").Value = 501
$ws.Range("D138").Value = 745490
$ws.Range("C141").Value = 10807
$ws.Range("D141").Value = 14410789
$ws.Range("C142").Value = 35033
$ws.Range("D142").Value = 50588990
$ws.Range("C143").Value = 81178
$ws.Range("D143").Value = 118935649
$ws.Range("C144").Value = 24327
$ws.Range("D144").Value = 36144275
$ws.Range("C145").Value = 6387
$ws.Range("D145").Value = 9530067
$ws.Range("C149").Value = 29188
$ws.Range("D149").Value = 39374480
